$d = $word.ActiveDocument
$d.Content.Find.Execute("80×53=4240", $true, $false, $false, $false, $false, $true, 1, $false, "43×88=3784", 2) | Out-Null
$d.Content.Find.Execute("88×43=3784", $true, $false, $false, $false, $false, $true, 1, $false, "60×73=4380", 2) | Out-Null
$d.Content.Find.Execute("36×95=3420", $true, $false, $false, $false, $false, $true, 1, $false, "96×41=3936", 2) | Out-Null
$d.Content.Find.Execute("37×61=2257", $true, $false, $false, $false, $false, $true, 1, $false, "18×49=882", 2) | Out-Null
$d.Content.Find.Execute("19×13=247", $true, $false, $false, $false, $false, $true, 1, $false, "85×70=5950", 2) | Out-Null
$d.Content.Find.Execute("70×99=6930", $true, $false, $false, $false, $false, $true, 1, $false, "66×77=5082", 2) | Out-Null
$d.Content.Find.Execute("59×83=4897", $true, $false, $false, $false, $false, $true, 1, $false, "86×62=5332", 2) | Out-Null
$d.Content.Find.Execute("26×19=494", $true, $false, $false, $false, $false, $true, 1, $false, "72×81=5832", 2) | Out-Null
$d.Content.Find.Execute("20×69=1380", $true, $false, $false, $false, $false, $true, 1, $false, "27×22=594", 2) | Out-Null
$d.Content.Find.Execute("26×89=2314", $true, $false, $false, $false, $false, $true, 1, $false, "74×84=6216", 2) | Out-Null
$d.Content.Find.Execute("25×58=1450", $true, $false, $false, $false, $false, $true, 1, $false, "65×66=4290", 2) | Out-Null
$d.Content.Find.Execute("25×81=2025", $true, $false, $false, $false, $false, $true, 1, $false, "18×45=810", 2) | Out-Null
$d.Content.Find.Execute("91×18=1638", $true, $false, $false, $false, $false, $true, 1, $false, "42×45=1890", 2) | Out-Null
$d.Content.Find.Execute("96×17=1632", $true, $false, $false, $false, $false, $true, 1, $false, "54×48=2592", 2) | Out-Null
$d.Content.Find.Execute("39×22=858", $true, $false, $false, $false, $false, $true, 1, $false, "63×92=5796", 2) | Out-Null
$d.Content.Find.Execute("32×78=2496", $true, $false, $false, $false, $false, $true, 1, $false, "86×73=6278", 2) | Out-Null
$d.Content.Find.Execute("45×56=2520", $true, $false, $false, $false, $false, $true, 1, $false, "47×37=1739", 2) | Out-Null
$d.Content.Find.Execute("24×69=1656", $true, $false, $false, $false, $false, $true, 1, $false, "46×35=1610", 2) | Out-Null
$d.Content.Find.Execute("72×29=2088", $true, $false, $false, $false, $false, $true, 1, $false, "85×13=1105", 2) | Out-Null
$d.Content.Find.Execute("15×81=1215", $true, $false, $false, $false, $false, $true, 1, $false, "66×53=3498", 2) | Out-Null
$d.Content.Find.Execute("25×99=2475", $true, $false, $false, $false, $false, $true, 1, $false, "25×15=375", 2) | Out-Null
$d.Content.Find.Execute("73×58=4234", $true, $false, $false, $false, $false, $true, 1, $false, "74×65=4810", 2) | Out-Null
$d.Content.Find.Execute("40×21=840", $true, $false, $false, $false, $false, $true, 1, $false, "36×81=2916", 2) | Out-Null
$d.Content.Find.Execute("56×28=1568", $true, $false, $false, $false, $false, $true, 1, $false, "28×75=2100", 2) | Out-Null
$d.Content.Find.Execute("87×11=957", $true, $false, $false, $false, $false, $true, 1, $false, "79×50=3950", 2) | Out-Null
